$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cells whose new text would otherwise be mis-parsed as a number:
# set NumberFormat to Text ("@") first, assign, then restore to the default
# "Normal" style so the saved file carries no stray style index.
$textCells = @('D5', 'D9', 'D10', 'D11', 'D15', 'D17', 'D20', 'D21', 'D22', 'D23', 'D25', 'D28', 'D29', 'D30', 'D34', 'D35', 'D37', 'D39', 'D43', 'D45', 'D46', 'D47', 'D48', 'D49', 'D50')
foreach ($addr in $textCells) { $ws.Range($addr).NumberFormat = "@" }
$ws.Range('D5').Value = '212.86'
$ws.Range('D9').Value = '0.0615'
$ws.Range('D10').Value = '18.41'
$ws.Range('D11').Value = '0.0815'
$ws.Range('D15').Value = '0.515'
$ws.Range('D17').Value = '61.83'
$ws.Range('D20').Value = '202.73'
$ws.Range('D21').Value = '4.29'
$ws.Range('D22').Value = '9.32'
$ws.Range('D23').Value = '6.02'
$ws.Range('D25').Value = '144.30'
$ws.Range('D28').Value = '15.21'
$ws.Range('D29').Value = '6.57'
$ws.Range('D30').Value = '0.0491'
$ws.Range('D34').Value = '2.43'
$ws.Range('D35').Value = '1.48'
$ws.Range('D37').Value = '0.0165'
$ws.Range('D39').Value = '0.791'
$ws.Range('D43').Value = '5.23'
$ws.Range('D45').Value = '91.70'
$ws.Range('D46').Value = '1.54'
$ws.Range('D47').Value = '54.28'
$ws.Range('D48').Value = '0.0507'
$ws.Range('D49').Value = '0.407'
$ws.Range('D50').Value = '1.00'
foreach ($addr in $textCells) { $ws.Range($addr).Style = "Normal" }

# --- Remaining cells: safe to assign directly (never parse as plain numbers).
$ws.Range('D2').Value = '26.296.22'
$ws.Range('E2').Value = '  +0.12%  '
$ws.Range('D3').Value = '1.611.01'
$ws.Range('E3').Value = '  +0.20%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('E6').Value = '  +0.02%  '
$ws.Range('E7').Value = '  -0.06%  '
$ws.Range('E8').Value = '  +0.01%  '
$ws.Range('E9').Value = '  -0.52%  '
$ws.Range('E10').Value = '  +1.85%  '
$ws.Range('E11').Value = '  -1.37%  '
$ws.Range('D12').Value = '1.836.59'
$ws.Range('E12').Value = '  +0.29%  '
$ws.Range('D13').Value = '1.609.69'
$ws.Range('E13').Value = '  +0.04%  '
$ws.Range('E14').Value = '  +0.36%  '
$ws.Range('E15').Value = '  +0.68%  '
$ws.Range('D16').Value = '26.297.70'
$ws.Range('E16').Value = '  +0.17%  '
$ws.Range('D18').Value = '0.0₃0729'
$ws.Range('E18').Value = '  +0.49%  '
$ws.Range('E19').Value = '  -0.04%  '
$ws.Range('E20').Value = '  +0.80%  '
$ws.Range('E21').Value = '  +0.65%  '
$ws.Range('E22').Value = '  -0.28%  '
$ws.Range('E23').Value = '  +0.18%  '
$ws.Range('E24').Value = '  +6.39%  '
$ws.Range('E25').Value = '  +1.42%  '
$ws.Range('E26').Value = '  -0.08%  '
$ws.Range('E27').Value = '  -2.06%  '
$ws.Range('E28').Value = '  -0.11%  '
$ws.Range('E29').Value = '  +1.59%  '
$ws.Range('E30').Value = '  +4.13%  '
$ws.Range('E31').Value = '  +0.38%  '
$ws.Range('E32').Value = '  +1.65%  '
$ws.Range('E33').Value = '  -1.67%  '
$ws.Range('E34').Value = '  +3.23%  '
$ws.Range('E35').Value = '  +0.42%  '
$ws.Range('D36').Value = '1.160.17'
$ws.Range('E36').Value = '  +5.03%  '
$ws.Range('E37').Value = '  +3.18%  '
$ws.Range('E38').Value = '  +0.06%  '
$ws.Range('E39').Value = '  +0.70%  '
$ws.Range('E40').Value = '  -0.32%  '
$ws.Range('E41').Value = '  +0.28%  '
$ws.Range('E42').Value = '  +0.92%  '
$ws.Range('E43').Value = '  +2.42%  '
$ws.Range('D44').Value = '1.750.48'
$ws.Range('E44').Value = '  +0.38%  '
$ws.Range('E45').Value = '  -1.54%  '
$ws.Range('E46').Value = '  -1.13%  '
$ws.Range('E47').Value = '  +1.24%  '
$ws.Range('E48').Value = '  +0.10%  '
$ws.Range('B49').Value = 'Mantle'
$ws.Range('C49').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('E49').Value = '  -0.20%  '
$ws.Range('B50').Value = 'USDD'
$ws.Range('C50').Value = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
$ws.Range('E50').Value = '  -0.13%  '
$ws.Range('B51').Value = 'BabyDogeCoin'
$ws.Range('C51').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D51').Value = '0.0₇0948'
$ws.Range('E51').Value = '  -15.58%  '
